$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Insert 4 new blank rows right after row 26 (before the old total
#    row, which was row 27) to make room for the new product lines.
# ------------------------------------------------------------------
$ws.Range("A27:A30").EntireRow.Insert()

# Clone the formatting (styles) of row 26 into the 4 freshly inserted
# rows so the new cells carry the same borders / number formats /
# fonts as the rest of the product table.
$ws.Range("A26:Q26").Copy()
$ws.Range("A27:Q30").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row heights follow the same alternating pattern as the existing
# product rows (25 -> 24.75, 26 -> 25.5, ...).
$ws.Rows.Item(27).RowHeight = 25.5
$ws.Rows.Item(28).RowHeight = 24.75
$ws.Rows.Item(29).RowHeight = 25.5
$ws.Rows.Item(30).RowHeight = 24.75

# Re-create the per-row merges (A:B, C:G, H:K, L:M, N:O) for the new
# rows, matching the pattern used by every other product row.
for ($r = 27; $r -le 30; $r++) {
    $ws.Range("A" + $r + ":B" + $r).Merge()
    $ws.Range("C" + $r + ":G" + $r).Merge()
    $ws.Range("H" + $r + ":K" + $r).Merge()
    $ws.Range("L" + $r + ":M" + $r).Merge()
    $ws.Range("N" + $r + ":O" + $r).Merge()
}

# ------------------------------------------------------------------
# 2. Update the two existing product rows whose item changed
#    (the new stock report re-sorted a couple of names into these
#    slots).
# ------------------------------------------------------------------
# Row 25: was "محلول ملح" -> now "صوفى طويل جدا جدا"
$ws.Range("C25").Value = "صوفى طويل جدا جدا"
$ws.Range("H25").Value = "7:0"
$ws.Range("L25").Value = "0"
$ws.Range("N25").Value = "55.00"
$ws.Range("P25").Value = "55.0000"
$ws.Range("Q25").Value = "1:0"

# Row 26: was "معجون سيجنال 50 مل" -> now "فرش اسنان اورل فريش"
$ws.Range("C26").Value = "فرش اسنان اورل فريش"
$ws.Range("H26").Value = "9:0"
$ws.Range("L26").Value = "0"
$ws.Range("N26").Value = "15.00"
$ws.Range("P26").Value = "15.0000"
$ws.Range("Q26").Value = "1:0"

# ------------------------------------------------------------------
# 3. Fill in the 4 new product rows.
# ------------------------------------------------------------------
# Row 27: ليفه
$ws.Range("A27").Value = 21
$ws.Range("C27").Value = "ليفه"
$ws.Range("H27").Value = "7:0"
$ws.Range("L27").Value = "0"
$ws.Range("N27").Value = "15.00"
$ws.Range("P27").Value = "15.0000"
$ws.Range("Q27").Value = "1:0"

# Row 28: محلول ملح
$ws.Range("A28").Value = 22
$ws.Range("C28").Value = "محلول ملح"
$ws.Range("H28").Value = "0:0"
$ws.Range("L28").Value = "0"
$ws.Range("N28").Value = "24.00"
$ws.Range("P28").Value = "24.0000"
$ws.Range("Q28").Value = "1:0"

# Row 29: مسواك اسنان (trailing space preserved, as in the source)
$ws.Range("A29").Value = 23
$ws.Range("C29").Value = "مسواك اسنان "
$ws.Range("H29").Value = "2:0"
$ws.Range("L29").Value = "0"
$ws.Range("N29").Value = "15.00"
$ws.Range("P29").Value = "15.0000"
$ws.Range("Q29").Value = "1:0"

# Row 30: معجون سيجنال 50 مل
$ws.Range("A30").Value = 24
$ws.Range("C30").Value = "معجون سيجنال 50 مل"
$ws.Range("H30").Value = "13:0"
$ws.Range("L30").Value = "0"
$ws.Range("N30").Value = "35.00"
$ws.Range("P30").Value = "35.0000"
$ws.Range("Q30").Value = "1:0"

# ------------------------------------------------------------------
# 4. Update the grand-total cell (now on row 31) to reflect the 4
#    newly added lines (+100.00).
# ------------------------------------------------------------------
$ws.Range("P31").Value = 1160.2550000000001

# ------------------------------------------------------------------
# 5. Update the footer timestamp (now on row 32) to the new export
#    time.
# ------------------------------------------------------------------
$ws.Range("A32").Value = "Tuesday, 26 August, 2025 12:09 PM"
